# Applies the scheduled-runner data update to the Leve profit tables
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
# Only raw data cells (H:N) are touched; no formulas exist in this workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 16666
$ws.Range("I21").Value = 10000
$ws.Range("K21").Value = 10000
$ws.Range("M21").Value = -9532
# Row 23
$ws.Range("H23").Value = 16666
$ws.Range("I23").Value = 10000
$ws.Range("K23").Value = 10000
$ws.Range("M23").Value = -9766
# Row 29
$ws.Range("H29").Value = 200
$ws.Range("I29").Value = 200
$ws.Range("K29").Value = 600
$ws.Range("M29").Value = -319
# Row 38
$ws.Range("H38").Value = 346.26666
$ws.Range("J38").Value = 854.5
$ws.Range("L38").Value = 2563.5
$ws.Range("N38").Value = -3307.5
# Row 58
$ws.Range("H58").Value = 1019.25
$ws.Range("I58").Value = 1019.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3057.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2907.75
$ws.Range("N58").ClearContents()
# Row 137
$ws.Range("H137").Value = 3624.1333
$ws.Range("I137").Value = 3447.8635
$ws.Range("J137").Value = 4108.875
$ws.Range("K137").Value = 10343.5905
$ws.Range("L137").Value = 12326.625
$ws.Range("M137").Value = -7793.5905
$ws.Range("N137").Value = -17426.625
# Row 141
$ws.Range("H141").Value = 5296.55
$ws.Range("I141").Value = 2988
$ws.Range("J141").Value = 7605.1
$ws.Range("K141").Value = 8964
$ws.Range("L141").Value = 22815.3
$ws.Range("M141").Value = -3784
$ws.Range("N141").Value = -33175.3

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
# Row 32
$ws.Range("H32").Value = 452111.66
$ws.Range("I32").Value = 582858.0600000001
$ws.Range("J32").Value = 10842.5625
$ws.Range("K32").Value = 582858.0600000001
$ws.Range("L32").Value = 10842.5625
$ws.Range("M32").Value = -582571.0600000001
$ws.Range("N32").Value = -11416.5625
# Row 37
$ws.Range("H37").Value = 11929.1
$ws.Range("J37").Value = 11929.1
$ws.Range("L37").Value = 11929.1
$ws.Range("N37").Value = -12475.1
# Row 44
$ws.Range("H44").Value = 17598.555
$ws.Range("J44").Value = 17598.555
$ws.Range("L44").Value = 17598.555
$ws.Range("N44").Value = -18574.555
# Row 55
$ws.Range("H55").Value = 27999.5
$ws.Range("J55").Value = 27999.5
$ws.Range("L55").Value = 27999.5
$ws.Range("N55").Value = -28629.5
# Row 61
$ws.Range("H61").Value = 2728.3845
$ws.Range("I61").Value = 2062.5
$ws.Range("J61").Value = 3793.8
$ws.Range("K61").Value = 2062.5
$ws.Range("L61").Value = 3793.8
$ws.Range("M61").Value = -1850.5
$ws.Range("N61").Value = -4217.8
# Row 74
$ws.Range("H74").Value = 1525.45
$ws.Range("I74").Value = 1471.2222
$ws.Range("J74").Value = 2013.5
$ws.Range("K74").Value = 1471.2222
$ws.Range("L74").Value = 2013.5
$ws.Range("M74").Value = -597.2221999999999
$ws.Range("N74").Value = -3761.5
# Row 77
$ws.Range("H77").Value = 1525.45
$ws.Range("I77").Value = 1471.2222
$ws.Range("J77").Value = 2013.5
$ws.Range("K77").Value = 7356.111
$ws.Range("L77").Value = 10067.5
$ws.Range("M77").Value = -2988.111
$ws.Range("N77").Value = -18803.5
# Row 135
$ws.Range("H135").Value = 48266
$ws.Range("J135").Value = 48266
$ws.Range("L135").Value = 48266
$ws.Range("N135").Value = -58406
# Row 136
$ws.Range("H136").Value = 2728.3845
$ws.Range("I136").Value = 2062.5
$ws.Range("J136").Value = 3793.8
$ws.Range("K136").Value = 6187.5
$ws.Range("L136").Value = 11381.4
$ws.Range("M136").Value = -3637.5
$ws.Range("N136").Value = -16481.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
# Row 134
$ws.Range("H134").Value = 4660.5713
$ws.Range("I134").Value = 4924.8
$ws.Range("K134").Value = 14774.4
$ws.Range("M134").Value = -12239.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
# Row 31
$ws.Range("H31").Value = 1676.6897
$ws.Range("I31").Value = 995.9545000000001
$ws.Range("K31").Value = 995.9545000000001
$ws.Range("M31").Value = -700.9545000000001
# Row 34
$ws.Range("H34").Value = 1676.6897
$ws.Range("I34").Value = 995.9545000000001
$ws.Range("K34").Value = 995.9545000000001
$ws.Range("M34").Value = -793.9545000000001
# Row 41
$ws.Range("H41").Value = 15899.223
$ws.Range("I41").Value = 100
$ws.Range("J41").Value = 17874.125
$ws.Range("K41").Value = 100
$ws.Range("L41").Value = 17874.125
$ws.Range("M41").Value = 328
$ws.Range("N41").Value = -18730.125
# Row 50
$ws.Range("H50").Value = 19249
$ws.Range("J50").Value = 19249
$ws.Range("L50").Value = 19249
$ws.Range("N50").Value = -20499
# Row 51
$ws.Range("H51").Value = 19249
$ws.Range("J51").Value = 19249
$ws.Range("L51").Value = 19249
$ws.Range("N51").Value = -20721
# Row 58
$ws.Range("H58").Value = 1357.0526
$ws.Range("I58").Value = 985
$ws.Range("J58").Value = 1528.7693
$ws.Range("K58").Value = 985
$ws.Range("L58").Value = 1528.7693
$ws.Range("M58").Value = -782
$ws.Range("N58").Value = -1934.7693
# Row 61
$ws.Range("H61").Value = 19249
$ws.Range("J61").Value = 19249
$ws.Range("L61").Value = 19249
$ws.Range("N61").Value = -19945
# Row 132
$ws.Range("H132").Value = 7577375.5
$ws.Range("I132").Value = 934.0714
$ws.Range("K132").Value = 2802.2142
$ws.Range("M132").Value = -272.2142000000003
# Row 134
$ws.Range("H134").Value = 2456.889
# Row 136
$ws.Range("H136").Value = 1357.0526
$ws.Range("I136").Value = 985
$ws.Range("J136").Value = 1528.7693
$ws.Range("K136").Value = 2955
$ws.Range("L136").Value = 4586.3079
$ws.Range("M136").Value = -405
$ws.Range("N136").Value = -9686.3079

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 13964.333
$ws.Range("J3").Value = 19931.5
$ws.Range("L3").Value = 59794.5
$ws.Range("N3").Value = -60018.5
# Row 122
$ws.Range("H122").Value = 5143.5
$ws.Range("I122").Value = 368.5
$ws.Range("K122").Value = 3316.5
$ws.Range("M122").Value = -866.5
# Row 124
$ws.Range("H124").Value = 2041.4166
$ws.Range("I124").Value = 1790
$ws.Range("J124").Value = 2064.2727
$ws.Range("K124").Value = 5370
$ws.Range("L124").Value = 6192.8181
$ws.Range("M124").Value = -460
$ws.Range("N124").Value = -16012.8181
# Row 125
$ws.Range("H125").Value = 2174.1428
$ws.Range("J125").Value = 2174.1428
$ws.Range("L125").Value = 6522.428400000001
$ws.Range("N125").Value = -16362.4284
# Row 138
$ws.Range("H138").Value = 3828.4614
$ws.Range("J138").Value = 4305.758
$ws.Range("L138").Value = 12917.274
$ws.Range("N138").Value = -23197.274
# Row 140
$ws.Range("H140").Value = 1977.5652
$ws.Range("I140").Value = 1236.5
$ws.Range("J140").Value = 3671.4285
$ws.Range("K140").Value = 3709.5
$ws.Range("L140").Value = 11014.2855
$ws.Range("M140").Value = 1470.5
$ws.Range("N140").Value = -21374.2855

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2595.6365
$ws.Range("I126").Value = 2420.6667
$ws.Range("J126").Value = 2805.6
$ws.Range("K126").Value = 7262.000100000001
$ws.Range("L126").Value = 8416.799999999999
$ws.Range("M126").Value = -4792.000100000001
$ws.Range("N126").Value = -13356.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3542.122
$ws.Range("I68").Value = 2299.5715
$ws.Range("J68").Value = 4846.8
$ws.Range("K68").Value = 2299.5715
$ws.Range("L68").Value = 4846.8
$ws.Range("M68").Value = -1550.5715
$ws.Range("N68").Value = -6344.8
# Row 71
$ws.Range("H71").Value = 3542.122
$ws.Range("I71").Value = 2299.5715
$ws.Range("J71").Value = 4846.8
$ws.Range("K71").Value = 11497.8575
$ws.Range("L71").Value = 24234
$ws.Range("M71").Value = -7753.8575
$ws.Range("N71").Value = -31722
# Row 132
$ws.Range("H132").Value = 3799.9395
$ws.Range("I132").Value = 3033.5557
$ws.Range("J132").Value = 4719.6
$ws.Range("K132").Value = 9100.667099999999
$ws.Range("L132").Value = 14158.8
$ws.Range("M132").Value = -6570.667099999999
$ws.Range("N132").Value = -19218.8
# Row 136
$ws.Range("H136").Value = 3097.7
$ws.Range("I136").Value = 4026.1667
$ws.Range("K136").Value = 12078.5001
$ws.Range("M136").Value = -9528.500100000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 4507036.5
$ws.Range("I132").Value = 2521.923
$ws.Range("K132").Value = 7565.768999999999
$ws.Range("M132").Value = -5035.768999999999
# Row 136
$ws.Range("H136").Value = 3496.5334
$ws.Range("I136").Value = 3191.6924
$ws.Range("K136").Value = 9575.0772
$ws.Range("M136").Value = -7025.0772

